$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 30; everything currently at/after row 30
# (rows 30-110) shifts down to rows 31-111.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 45133
$ws.Range("D30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108007
$ws.Range("J30").Value = "Coco"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 36000
$ws.Range("O30").Value = 38000
$ws.Range("P30").Value = 36667
$ws.Range("Q30").Value = "$/malla 20 unidades"
$ws.Range("R30").Value = "Perú"
$ws.Range("S30").Value = 1833
$ws.Range("T30").Value = 20
